{"js": "// Melanjutkan pembuatan model sppd 1,2\n//\n// 1) The table's 7th row (index 6, 0-based) has a blank first cell that\n//    should read \"7\" (continuing the numbered rows 1..6, 8, 9, 10...).\n// 2) Five new paragraphs are appended after the table (right before the\n//    section break) describing where/when the document was issued and\n//    who signed it.\n\nconst body = context.document.body;\n\n// --- 1) Fill in the missing row number \"7\" ---------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst targetCell = table.getCell(6, 0);\ntargetCell.value = \"7\";\n\n// --- 2) Append the closing/signature paragraphs after the table ------------\nbody.insertParagraph(\"Dikeluarkan di : Palembang\", \"End\");\nbody.insertParagraph(\"Pada tanggal   : September 2022\", \"End\");\nbody.insertParagraph(\"Sekretaris DPRD Provinsi Sumatera Selatan\", \"End\");\nbody.insertParagraph(\"Ramdhan s User\", \"End\");\nbody.insertParagraph(\"Pembina Utama Madya (NIP : 090419204109409)\", \"End\");\n\nawait context.sync();\n", "ps1": "# Melanjutkan pembuatan model sppd 1,2\n#\n# 1) The table's 7th row has a blank first-column cell that should read\n#    \"7\" (continuing the numbered rows 1..6, 8, 9, 10...).\n# 2) Five new paragraphs are appended after the table (right before the\n#    section break) describing where/when the document was issued and\n#    who signed it.\n\n$d = $word.ActiveDocument\n\n# --- 1) Fill in the missing row number \"7\" ---------------------------------\n$table = $d.Tables.Item(1)\n$cell = $table.Cell(7, 1)\n$cell.Range.Text = \"7\"\n\n# --- 2) Append the closing/signature paragraphs after the table ------------\n$lines = @(\n  \"Dikeluarkan di : Palembang\",\n  \"Pada tanggal   : September 2022\",\n  \"Sekretaris DPRD Provinsi Sumatera Selatan\",\n  \"Ramdhan s User\",\n  \"Pembina Utama Madya (NIP : 090419204109409)\"\n)\n\nforeach ($line in $lines) {\n  $end = $d.Content\n  $end.Collapse(0)          # wdCollapseEnd\n  $end.InsertParagraphAfter()\n  $end.Collapse(0)          # wdCollapseEnd\n  $end.Text = $line\n}\n"}
